$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Unprotect()

# Insert a new row above row 7 (shifts header/data rows down by one)
$ws.Rows.Item(7).Insert()

# Populate the new row 7 with the "Help" label and hyperlink text
$ws.Range("A7").Value = "Help"
$ws.Hyperlinks.Add($ws.Range("B7"), "https://ccdb.esss.lu.se/resources/help/ccdb_conventions.pdf", "", "", "Help")
